$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing "Projects Final Score" column (F),
# shifting it to G, and use the freed-up F column for a new "Final Exam"
# column.
$ws.Columns("F:F").Insert()

# Header for the new column.
$ws.Range("F1").Value = " Final Exam"

# New "Final Exam" scores for the 5 students (rows 2-6).
$ws.Range("F2").Value = 80
$ws.Range("F3").Value = 65
$ws.Range("F4").Value = 92
$ws.Range("F5").Value = 78
$ws.Range("F6").Value = 63

# Column width tweaks (set widths first, then hide column F - changing
# order avoids emitting duplicate <col> entries for the same column).
$ws.Columns("A:A").ColumnWidth = 24.333333333333332
$ws.Columns("B:B").ColumnWidth = 13.5
$ws.Columns("E:E").ColumnWidth = 19.666666666666668
$ws.Columns("F:F").ColumnWidth = 14.833333333333334
$ws.Columns("G:G").ColumnWidth = 12.5

# Hide the new column (matches the other raw-score columns C:E which are
# hidden helper columns).
$ws.Columns("F:F").Hidden = $true

# A handful of leftover/paste-artifact rows that appeared beneath the first
# small table.
$ws.Range("B16").Font.Name = "Arial"
$ws.Range("B16").Font.Size = 10
$ws.Range("C16").Font.Name = "Arial"
$ws.Range("C16").Font.Size = 10
$ws.Range("C16").Value = "Projects Final Score"

$ws.Range("B17").Font.Name = "Arial"
$ws.Range("B17").Font.Size = 10
$ws.Range("C17").Font.Name = "Arial"
$ws.Range("C17").Font.Size = 10
$ws.Range("C17").Value = "Projects Final Score"

$ws.Range("B18").Font.Name = "Arial"
$ws.Range("B18").Font.Size = 10
$ws.Range("C18").Font.Name = "Arial"
$ws.Range("C18").Font.Size = 10
$ws.Range("C18").Value = "Projects Final Score"

$ws.Range("B19").Font.Name = "Arial"
$ws.Range("B19").Font.Size = 10
$ws.Range("C19").Font.Name = "Arial"
$ws.Range("C19").Font.Size = 10
$ws.Range("C19").Value = "Projects Final Score"

# Selection ends on G1, matching the saved workbook view.
$ws.Range("G1").Select()
